# Arbeitszeit_Pichler.xlsx update:
# - Append 5 new time-tracking rows (24-28) covering 2019-08-02..2019-08-06
# - New Tätigkeit/Kommentar text values land as new shared strings
# - B6 (SUM(H:H)) / B7 recalc automatically from the new H column entries
# - Update the active selection to J29 (matches the author's last click)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: 02.08.2019, 2h Programmieren
$ws.Range("G23").Copy($ws.Range("G24"))
$ws.Range("G24").Value = 43679
$ws.Range("H24").Value = 2
$ws.Range("I24").Value = "Stunden"
$ws.Range("J24").Value = "Programmieren"
$ws.Range("K24").Value = "Schwierigkeiten bei der Implementierung der Bt Bibliothek"

# Row 25: 03.08.2019, 1h Programmieren
$ws.Range("G23").Copy($ws.Range("G25"))
$ws.Range("G25").Value = 43680
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = "Stunden"
$ws.Range("J25").Value = "Programmieren"

# Row 26: 04.08.2019, 2h Dokumentation
$ws.Range("G23").Copy($ws.Range("G26"))
$ws.Range("G26").Value = 43681
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = "Stunden"
$ws.Range("J26").Value = "Dokumentation"
$ws.Range("K26").Value = "Hauptsächlich Dokumentation der vorhergegangenen Programmiertätigkeit"

# Row 27: 05.08.2019, 1h Programmieren
$ws.Range("G23").Copy($ws.Range("G27"))
$ws.Range("G27").Value = 43682
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = "Stunden"
$ws.Range("J27").Value = "Programmieren"
$ws.Range("K27").Value = "Erweiterung der StreamClient Klasse"

# Row 28: 06.08.2019, 1h Programmieren
$ws.Range("G23").Copy($ws.Range("G28"))
$ws.Range("G28").Value = 43683
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = "Stunden"
$ws.Range("J28").Value = "Programmieren"
$ws.Range("K28").Value = "Erstellen eines PS Scripts, siehe vorheriger Eintrag"

# Scroll/selection state to match the author's saved view
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J29").Select()
